# -----------------------------------------------------------------------
# "Ran code for averaged intensites on spiral schemes"
#
# The underlying generation script was re-run after adding three new
# "Spiral" orientation schemes (Gaussian-Quadrature moves up next to the
# other single-point schemes, the three new Spiral-* rows are inserted
# right after it, and all the previously-existing rotation/hex-grid rows
# shift down by four rows). The sheet grows from 16 to 19 used rows
# (A1:P16 -> A1:P19).
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the three brand-new HKL-index cells (A17:A19) the same bold /
# centered / bordered formatting already used by the rest of column A
# (style of the last existing index cell, A16) before the values land.
$ws.Range("A16").Copy($ws.Range("A17"))
$ws.Range("A16").Copy($ws.Range("A18"))
$ws.Range("A16").Copy($ws.Range("A19"))
$wb.Application.CutCopyMode = $false

# Row layout (rows 10-19), each row: HKL index (col A), scheme name
# (col B) and the 14 averaged-intensity ratios (cols C:P).
$rowsData = @(
    @{ Row = 10; A = 8; Label = "Gaussian-Quadrature"; Vals = @(1.001045738816186, 0.9763106660148633, 0.9905882352941177, 0.9888527775402498, 1.001045738816186, 0.9763106660148633, 0.99, 0.981764705882353, 0.9925490388345444, 0.9819315554119415, 1.001045738816186, 0.9834494506544904, 0.9891993544163542, 0.9878803397242819) }
    @{ Row = 11; A = 9; Label = "Spiral-90deg-10rot-5space"; Vals = @(1.001267848907582, 0.9654873432150899, 0.9925839141481579, 0.9848410030365421, 1.001267848907582, 0.9654873432150899, 0.9931767631831939, 0.9839901330880197, 0.9965198442985799, 0.9735804460606551, 1.001267848907582, 0.979035628681624, 0.9860450273268431, 0.9864309119922277) }
    @{ Row = 12; A = 10; Label = "Spiral-90deg-15rot-5space"; Vals = @(1.001223648239131, 0.9655811943862864, 0.9925405557528618, 0.9848594304005645, 1.001223648239131, 0.9655811943862864, 0.9931308652924401, 0.9840247282631469, 0.9965100009819182, 0.9736416568864202, 1.001223648239131, 0.9790608750695742, 0.986051207194711, 0.9864390100253462) }
    @{ Row = 13; A = 11; Label = "Spiral-90deg-10rot-3space"; Vals = @(1.001264175234741, 0.9654797020836825, 0.9925779251043916, 0.984850326492912, 1.001264175234741, 0.9654797020836825, 0.993167078723062, 0.9840209221568416, 0.9965178978921707, 0.9735910168001257, 1.001264175234741, 0.979028813594037, 0.9860430322289319, 0.9864336305609909) }
    @{ Row = 14; A = 12; Label = "NoRotation-tilt60deg"; Vals = @(1, 0.9622840000000006, 1.001256000000001, 0.9785320000000013, 1, 0.9622840000000006, 0.9987600000000005, 0.9770800000000008, 0.9981240000000015, 0.9695760000000004, 1, 0.981770000000001, 0.9855180000000008, 0.9857015000000008) }
    @{ Row = 15; A = 13; Label = "Rotation-NoTilt"; Vals = @(1, 0.9438874999999985, 1.01, 0.97, 1, 0.9438874999999985, 1.01, 0.97, 1, 0.96, 1, 0.9769437499999993, 0.9809718749999996, 0.9829859374999999) }
    @{ Row = 16; A = 14; Label = "Rotation-60detTilt"; Vals = @(0.995259822899203, 0.9623750508544001, 1.000730705510398, 0.9781092528128005, 0.995259822899203, 0.9623750508544001, 1.000701407231998, 0.9781973538816008, 0.9945386862592011, 0.9716225710080002, 0.995259822899203, 0.981552878182399, 0.9841187080192003, 0.9851918563072002) }
    @{ Row = 17; A = 15; Label = "HexGrid-90degTilt5degRes"; Vals = @(0.987916960427406, 0.9882814289403906, 0.9877268623732334, 0.9880020122975326, 0.987916960427406, 0.9882814289403906, 0.9879949185940951, 0.9883449923671825, 0.9876332335055968, 0.987594492097427, 0.987916960427406, 0.988004145656812, 0.9879818160096406, 0.987936862575358) }
    @{ Row = 18; A = 16; Label = "HexGrid-90degTilt22p5degRes"; Vals = @(0.9874915922491165, 0.9896617477086177, 0.9875652911014803, 0.9877906375276257, 0.9874915922491165, 0.9896617477086177, 0.9881038456802084, 0.9880218876675262, 0.9864377240242617, 0.9889037895740711, 0.9874915922491165, 0.9886135194050489, 0.9881273171467101, 0.9879970644416134) }
    @{ Row = 19; A = 17; Label = "HexGrid-60degTilt5degRes"; Vals = @(0.9850077410494816, 0.9938160177949179, 0.9863556726472777, 0.9889384290707812, 0.9850077410494816, 0.9938160177949179, 0.9860485660581, 0.9890761757156208, 0.9864491816569338, 0.9910889678263319, 0.9850077410494816, 0.9900858452210979, 0.9885294651406146, 0.9883475939774307) }
)

foreach ($r in $rowsData) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.A
    $ws.Cells.Item($r.Row, 2).Value2 = $r.Label
    $col = 3
    foreach ($v in $r.Vals) {
        $ws.Cells.Item($r.Row, $col).Value2 = $v
        $col++
    }
}
